$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.618.36'
$ws.Range("E2").Value = '  -2.75%  '
$ws.Range("D3").Value = '2.397.38'
$ws.Range("E3").Value = '  -1.97%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = "'569.32"
$ws.Range("E5").Value = '  -2.36%  '
$ws.Range("D6").Value = "'139.69"
$ws.Range("E6").Value = '  -3.12%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("E8").Value = '  -1.56%  '
$ws.Range("D9").Value = '2.376.83'
$ws.Range("E9").Value = '  -2.72%  '
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = "'5.04"
$ws.Range("E12").Value = '  -3.34%  '
$ws.Range("E13").Value = '  -2.69%  '
$ws.Range("D14").Value = "'25.89"
$ws.Range("E14").Value = '  -2.55%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = "'0.0000169"
$ws.Range("E15").Value = '  -2.19%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.828.15'
$ws.Range("E16").Value = '  -2.05%  '
$ws.Range("D17").Value = '60.637.72'
$ws.Range("E17").Value = '  -2.52%  '
$ws.Range("D18").Value = '2.397.17'
$ws.Range("E18").Value = '  -1.81%  '
$ws.Range("D19").Value = "'10.51"
$ws.Range("E19").Value = '  -3.69%  '
$ws.Range("D20").Value = "'7.11"
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("D21").Value = "'320.47"
$ws.Range("E21").Value = '  -2.97%  '
$ws.Range("D22").Value = "'4.00"
$ws.Range("E22").Value = '  -2.94%  '
$ws.Range("D23").Value = "'6.09"
$ws.Range("E23").Value = '  +1.55%  '
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("D25").Value = "'1.87"
$ws.Range("E25").Value = '  -6.28%  '
$ws.Range("D26").Value = "'64.43"
$ws.Range("E26").Value = '  -2.31%  '
$ws.Range("E27").Value = '  -8.23%  '
$ws.Range("D28").Value = "'574.93"
$ws.Range("E28").Value = '  -8.05%  '
$ws.Range("D29").Value = '2.514.92'
$ws.Range("E29").Value = '  -1.93%  '
$ws.Range("D30").Value = '0.0₃0902'
$ws.Range("E30").Value = '  -6.03%  '
$ws.Range("D31").Value = "'7.78"
$ws.Range("E31").Value = '  -3.03%  '
$ws.Range("E32").Value = '  -7.42%  '
$ws.Range("E33").Value = '  -3.21%  '
$ws.Range("D34").Value = "'0.132"
$ws.Range("E34").Value = '  -7.42%  '
$ws.Range("D36").Value = "'4.59"
$ws.Range("E36").Value = '  -6.89%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = "'0.363"
$ws.Range("E37").Value = '  -3.82%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = "'1.37"
$ws.Range("E38").Value = '  -4.70%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = "'147.45"
$ws.Range("E39").Value = '  -1.58%  '
$ws.Range("D40").Value = "'18.06"
$ws.Range("E40").Value = '  -1.64%  '
$ws.Range("D41").Value = "'5.06"
$ws.Range("E41").Value = '  -4.94%  '
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("D43").Value = "'41.49"
$ws.Range("E43").Value = '  -2.31%  '
$ws.Range("D44").Value = "'1.65"
$ws.Range("E44").Value = '  -5.92%  '
$ws.Range("D45").Value = "'2.34"
$ws.Range("E45").Value = '  -5.50%  '
$ws.Range("D46").Value = '0.0₆0288'
$ws.Range("E46").Value = '  +20.37%  '
$ws.Range("D47").Value = "'140.01"
$ws.Range("E47").Value = '  -2.60%  '
$ws.Range("D48").Value = "'3.48"
$ws.Range("E48").Value = '  -4.40%  '
$ws.Range("D49").Value = "'0.583"
$ws.Range("E49").Value = '  -3.07%  '
$ws.Range("E50").Value = '  -4.44%  '
$ws.Range("D51").Value = "'19.27"
$ws.Range("E51").Value = '  -1.55%  '
